$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1700.091
$ws.Range("I98").Value = 1576.9231
$ws.Range("J98").Value = 1878
$ws.Range("K98").Value = 1576.9231
$ws.Range("L98").Value = 1878
$ws.Range("M98").Value = -78.92309999999998
$ws.Range("N98").Value = -4874
$ws.Range("H122").Value = 1700.091
$ws.Range("I122").Value = 1576.9231
$ws.Range("J122").Value = 1878
$ws.Range("K122").Value = 4730.7693
$ws.Range("L122").Value = 5634
$ws.Range("M122").Value = -2280.7693
$ws.Range("N122").Value = -10534

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 36.333332
$ws.Range("I4").Value = 29.5
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 29.5
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = 86.5
$ws.Range("N4").Value = -282
$ws.Range("H32").Value = 5761219
$ws.Range("I32").Value = 8490
$ws.Range("K32").Value = 8490
$ws.Range("M32").Value = -8203
$ws.Range("H45").Value = 2592.75
$ws.Range("I45").Value = 1422.2778
$ws.Range("J45").Value = 4097.643
$ws.Range("K45").Value = 1422.2778
$ws.Range("L45").Value = 4097.643
$ws.Range("M45").Value = -1045.2778
$ws.Range("N45").Value = -4851.643
$ws.Range("H61").Value = 2513.2703
$ws.Range("I61").Value = 1458.1904
$ws.Range("J61").Value = 3898.0625
$ws.Range("K61").Value = 1458.1904
$ws.Range("L61").Value = 3898.0625
$ws.Range("M61").Value = -1246.1904
$ws.Range("N61").Value = -4322.0625
$ws.Range("H74").Value = 71432830
$ws.Range("I74").Value = 107147890
$ws.Range("K74").Value = 107147890
$ws.Range("M74").Value = -107147016
$ws.Range("H77").Value = 71432830
$ws.Range("I77").Value = 107147890
$ws.Range("K77").Value = 535739450
$ws.Range("M77").Value = -535735082
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
$ws.Range("H132").Value = 1987586.2
$ws.Range("I132").Value = 1661
$ws.Range("J132").Value = 4469993
$ws.Range("K132").Value = 4983
$ws.Range("L132").Value = 13409979
$ws.Range("M132").Value = -2453
$ws.Range("N132").Value = -13415039
$ws.Range("H136").Value = 2513.2703
$ws.Range("I136").Value = 1458.1904
$ws.Range("J136").Value = 3898.0625
$ws.Range("K136").Value = 4374.5712
$ws.Range("L136").Value = 11694.1875
$ws.Range("M136").Value = -1824.5712
$ws.Range("N136").Value = -16794.1875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 544.069
$ws.Range("I22").Value = 603.25
$ws.Range("J22").Value = 260
$ws.Range("K22").Value = 603.25
$ws.Range("L22").Value = 260
$ws.Range("M22").Value = -430.25
$ws.Range("N22").Value = -606
$ws.Range("H134").Value = 4510.268
$ws.Range("I134").Value = 2081.7585
$ws.Range("J134").Value = 7118.6665
$ws.Range("K134").Value = 6245.2755
$ws.Range("L134").Value = 21355.9995
$ws.Range("M134").Value = -3710.2755
$ws.Range("N134").Value = -26425.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7144886
$ws.Range("I31").Value = 1329.7561
$ws.Range("J31").Value = 17244398
$ws.Range("K31").Value = 1329.7561
$ws.Range("L31").Value = 17244398
$ws.Range("M31").Value = -1034.7561
$ws.Range("N31").Value = -17244988
$ws.Range("H34").Value = 7144886
$ws.Range("I34").Value = 1329.7561
$ws.Range("J34").Value = 17244398
$ws.Range("K34").Value = 1329.7561
$ws.Range("L34").Value = 17244398
$ws.Range("M34").Value = -1127.7561
$ws.Range("N34").Value = -17244802
$ws.Range("H58").Value = 2391888.8
$ws.Range("I58").Value = 2908
$ws.Range("J58").Value = 3347481
$ws.Range("K58").Value = 2908
$ws.Range("L58").Value = 3347481
$ws.Range("M58").Value = -2705
$ws.Range("N58").Value = -3347887
$ws.Range("H122").Value = 66669770
$ws.Range("I122").Value = 83334540
$ws.Range("J122").Value = 10666.667
$ws.Range("K122").Value = 250003620
$ws.Range("L122").Value = 32000.001
$ws.Range("M122").Value = -250001170
$ws.Range("N122").Value = -36900.001
$ws.Range("H132").Value = 3466.2104
$ws.Range("I132").Value = 3172
$ws.Range("J132").Value = 3571.2856
$ws.Range("K132").Value = 9516
$ws.Range("L132").Value = 10713.8568
$ws.Range("M132").Value = -6986
$ws.Range("N132").Value = -15773.8568
$ws.Range("H134").Value = 2823.9656
$ws.Range("I134").Value = 1461.4615
$ws.Range("K134").Value = 4384.3845
$ws.Range("M134").Value = -1849.3845
$ws.Range("H136").Value = 2391888.8
$ws.Range("I136").Value = 2908
$ws.Range("J136").Value = 3347481
$ws.Range("K136").Value = 8724
$ws.Range("L136").Value = 10042443
$ws.Range("M136").Value = -6174
$ws.Range("N136").Value = -10047543

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 32
$ws.Range("I2").Value = 46.2
$ws.Range("K2").Value = 46.2
$ws.Range("M2").Value = 66.8
$ws.Range("H11").Value = 3815200
$ws.Range("I11").Value = 2351266.8
$ws.Range("J11").Value = 8207000
$ws.Range("K11").Value = 2351266.8
$ws.Range("L11").Value = 8207000
$ws.Range("M11").Value = -2351127.8
$ws.Range("N11").Value = -8207278
$ws.Range("H18").Value = 70006
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 70006
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 70006
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -70592
$ws.Range("H102").Value = 1833734.8
$ws.Range("I102").Value = 2858767.8
$ws.Range("J102").Value = 3318.6428
$ws.Range("K102").Value = 2858767.8
$ws.Range("L102").Value = 3318.6428
$ws.Range("M102").Value = -2857145.8
$ws.Range("N102").Value = -6562.6428

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 680.9091
$ws.Range("I9").Value = 244.28572
$ws.Range("J9").Value = 1445
$ws.Range("K9").Value = 244.28572
$ws.Range("L9").Value = 1445
$ws.Range("M9").Value = -20.28572
$ws.Range("N9").Value = -1893
$ws.Range("H13").Value = 7120.5
$ws.Range("I13").Value = 2222.2222
$ws.Range("J13").Value = 15937.4
$ws.Range("K13").Value = 2222.2222
$ws.Range("L13").Value = 15937.4
$ws.Range("M13").Value = -2082.2222
$ws.Range("N13").Value = -16217.4
$ws.Range("H22").Value = 1012.619
$ws.Range("I22").Value = 878.4545000000001
$ws.Range("J22").Value = 1160.2
$ws.Range("K22").Value = 878.4545000000001
$ws.Range("L22").Value = 1160.2
$ws.Range("M22").Value = -583.4545000000001
$ws.Range("N22").Value = -1750.2
$ws.Range("H27").Value = 1012.619
$ws.Range("I27").Value = 878.4545000000001
$ws.Range("J27").Value = 1160.2
$ws.Range("K27").Value = 878.4545000000001
$ws.Range("L27").Value = 1160.2
$ws.Range("M27").Value = -771.4545000000001
$ws.Range("N27").Value = -1374.2
$ws.Range("H40").Value = 63127560
$ws.Range("I40").Value = 77693550
$ws.Range("J40").Value = 8268.333000000001
$ws.Range("K40").Value = 77693550
$ws.Range("L40").Value = 8268.333000000001
$ws.Range("M40").Value = -77693414
$ws.Range("N40").Value = -8540.333000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 46716.31
$ws.Range("J46").Value = 46716.31
$ws.Range("L46").Value = 46716.31
$ws.Range("N46").Value = -47178.31
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 10000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 10000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9685
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 10000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 10000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8908
$ws.Range("N73").ClearContents()
$ws.Range("H122").Value = 6113.125
$ws.Range("I122").Value = 6505.2
$ws.Range("K122").Value = 19515.6
$ws.Range("M122").Value = -17065.6
$ws.Range("H126").Value = 4992.2573
$ws.Range("I126").Value = 5208.5
$ws.Range("J126").Value = 2685.6667
$ws.Range("K126").Value = 15625.5
$ws.Range("L126").Value = 8057.000100000001
$ws.Range("M126").Value = -13155.5
$ws.Range("N126").Value = -12997.0001
$ws.Range("H132").Value = 2094.5107
$ws.Range("I132").Value = 1358.3334
$ws.Range("J132").Value = 3393.647
$ws.Range("K132").Value = 4075.0002
$ws.Range("L132").Value = 10180.941
$ws.Range("M132").Value = -1545.0002
$ws.Range("N132").Value = -15240.941
$ws.Range("H134").Value = 46716.31
$ws.Range("J134").Value = 46716.31
$ws.Range("L134").Value = 140148.93
$ws.Range("N134").Value = -145218.93
$ws.Range("H136").Value = 9269755
$ws.Range("I136").Value = 10427870
$ws.Range("J136").Value = 4831.6665
$ws.Range("K136").Value = 31283610
$ws.Range("L136").Value = 14494.9995
$ws.Range("M136").Value = -31281060
$ws.Range("N136").Value = -19594.9995
